# Added New Mac-Address and Document Types
# Appends a new test-data row (row 33) to the master-reg_center_user sheet,
# mirroring the existing rows' layout (regcntr_id, usr_id, lang_code,
# is_active, cr_by, cr_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item(33, 1).Row

$ws.Cells.Item($newRow, 1).Value = 10002
$ws.Cells.Item($newRow, 2).Value = 110032
$ws.Cells.Item($newRow, 3).Value = "eng"
$ws.Cells.Item($newRow, 4).Value = $true
$ws.Cells.Item($newRow, 5).Value = "superadmin"
$ws.Cells.Item($newRow, 6).Value = "now()"

# Reflect the post-edit selection/scroll position seen after adding the row
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 22
